$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = @{ D = 751258.5067377688;  E = 132643.9025941133 }
    3  = @{ D = 141224.85510650915; E = 13848.317991206073 }
    4  = @{ D = 71016.18235385128;  E = 9284.533410347995 }
    5  = @{ D = 539017.4692774098;  E = 109511.0511925592 }
    6  = @{ D = 27578.037147134797; E = 5740.358000462053 }
    7  = @{ D = 133224.93828604417; E = 24807.243343738122 }
    8  = @{ D = 189511.85823830997; E = 32352.087319831586 }
    9  = @{ D = 222768.88532188482; E = 37595.732117752814 }
    10 = @{ D = 178028.82456520526; E = 32120.997452340147 }
    11 = @{ D = 145.96317919143652; E = 27.484359988560563 }
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row].D
    $ws.Range("E$row").Value = $values[$row].E
}
